$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 0.3876063333333333
$ws.Range("H2").Value = 1.162819
$ws.Range("I2").Value = 0.009710908683933402
$ws.Range("J2").Value = 0.009710908683933402
$ws.Range("M2").Value = 8.033114333333334
$ws.Range("N2").Value = 24.099343
$ws.Range("O2").Value = 0.1374088679258946
$ws.Range("P2").Value = 0.1374088679258946
$ws.Range("Q2").Value = 3.113685991990778
$ws.Range("R2").Value = 28.023173927917
$ws.Range("S2").Value = 0.001334364968791027
$ws.Range("T2").Value = 0.001334364968791028
$ws.Range("G3").Value = 0.3876063333333333
$ws.Range("H3").Value = 1.162819
$ws.Range("I3").Value = 0.009710908683933402
$ws.Range("J3").Value = 0.009710908683933402
$ws.Range("O3").Value = 0.6355200716780686
$ws.Range("P3").Value = 0.6355200716780686
$ws.Range("Q3").Value = 14.40088965640967
$ws.Range("R3").Value = 129.608006907687
$ws.Range("S3").Value = 0.006171477382872534
$ws.Range("T3").Value = 0.006171477382872534
$ws.Range("G4").Value = 0.3876063333333333
$ws.Range("H4").Value = 1.162819
$ws.Range("I4").Value = 0.009710908683933402
$ws.Range("J4").Value = 0.009710908683933402
$ws.Range("M4").Value = 13.27489133333333
$ws.Range("N4").Value = 39.824674
$ws.Range("O4").Value = 0.2270710603960369
$ws.Range("P4").Value = 0.2270710603960369
$ws.Range("Q4").Value = 5.145431955111778
$ws.Range("R4").Value = 46.308887596006
$ws.Range("S4").Value = 0.00220506633226984
$ws.Range("T4").Value = 0.00220506633226984
$ws.Range("I5").Value = 0.714669937678414
$ws.Range("J5").Value = 0.714669937678414
$ws.Range("M5").Value = 8.033114333333334
$ws.Range("N5").Value = 24.099343
$ws.Range("O5").Value = 0.1374088679258946
$ws.Range("P5").Value = 0.1374088679258946
$ws.Range("Q5").Value = 229.150313968853
$ws.Range("R5").Value = 2062.352825719677
$ws.Range("S5").Value = 0.0982019870770605
$ws.Range("T5").Value = 0.09820198707706052
$ws.Range("I6").Value = 0.714669937678414
$ws.Range("J6").Value = 0.714669937678414
$ws.Range("O6").Value = 0.6355200716780686
$ws.Range("P6").Value = 0.6355200716780686
$ws.Range("S6").Value = 0.4541870900195465
$ws.Range("T6").Value = 0.4541870900195465
$ws.Range("I7").Value = 0.714669937678414
$ws.Range("J7").Value = 0.714669937678414
$ws.Range("M7").Value = 13.27489133333333
$ws.Range("N7").Value = 39.824674
$ws.Range("O7").Value = 0.2270710603960369
$ws.Range("P7").Value = 0.2270710603960369
$ws.Range("Q7").Value = 378.675740280854
$ws.Range("R7").Value = 3408.081662527686
$ws.Range("S7").Value = 0.162280860581807
$ws.Range("T7").Value = 0.162280860581807
$ws.Range("G8").Value = 11.00120833333333
$ws.Range("H8").Value = 33.003625
$ws.Range("I8").Value = 0.2756191536376525
$ws.Range("J8").Value = 0.2756191536376525
$ws.Range("M8").Value = 8.033114333333334
$ws.Range("N8").Value = 24.099343
$ws.Range("O8").Value = 0.1374088679258946
$ws.Range("P8").Value = 0.1374088679258946
$ws.Range("Q8").Value = 88.3739643464861
$ws.Range("R8").Value = 795.3656791183751
$ws.Range("S8").Value = 0.03787251588004304
$ws.Range("T8").Value = 0.03787251588004305
$ws.Range("G9").Value = 11.00120833333333
$ws.Range("H9").Value = 33.003625
$ws.Range("I9").Value = 0.2756191536376525
$ws.Range("J9").Value = 0.2756191536376525
$ws.Range("O9").Value = 0.6355200716780686
$ws.Range("P9").Value = 0.6355200716780686
$ws.Range("Q9").Value = 408.7321946807916
$ws.Range("R9").Value = 3678.589752127125
$ws.Range("S9").Value = 0.1751615042756495
$ws.Range("T9").Value = 0.1751615042756495
$ws.Range("G10").Value = 11.00120833333333
$ws.Range("H10").Value = 33.003625
$ws.Range("I10").Value = 0.2756191536376525
$ws.Range("J10").Value = 0.2756191536376525
$ws.Range("M10").Value = 13.27489133333333
$ws.Range("N10").Value = 39.824674
$ws.Range("O10").Value = 0.2270710603960369
$ws.Range("P10").Value = 0.2270710603960369
$ws.Range("Q10").Value = 146.0398451603611
$ws.Range("R10").Value = 1314.35860644325
$ws.Range("S10").Value = 0.06258513348195996
$ws.Range("T10").Value = 0.06258513348195996

Write-Output "Applied updated TPM values"
